# Update cryptos list data (price + 1h volume%) and fix two swapped rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.189.64'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.764.85'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '354.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.95%  '
$ws.Range('E7').Value = '  -2.69%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.43'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('E11').Value = '  +3.57%  '
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0834'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.94'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.52'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.199.81'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.755.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.929'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.140.99'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.08'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.00%  '
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.64'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '264.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.71%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.99%  '
$ws.Range('E28').Value = '  +13.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.68%  '
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '51.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '34.74'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.06'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.59%  '
$ws.Range('E34').Value = '  -3.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.54'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.65%  '
$ws.Range('E36').Value = '  -0.51%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.17'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('E39').Value = '  -1.85%  '
$ws.Range('E40').Value = '  -1.75%  '
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '120.37'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.30%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.05'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.97%  '
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.085.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.25%  '
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('E48').Value = '  -0.86%  '
$ws.Range('E49').Value = '  -3.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.915'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('E51').Value = '  +6.51%  '
